# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue ($ws.Range('D2')) '28.150.95'
Set-TextValue ($ws.Range('E2')) '  +3.47%  '
Set-TextValue ($ws.Range('D3')) '1.577.28'
Set-TextValue ($ws.Range('E4')) '  -0.94%  '
Set-TextValue ($ws.Range('D5')) '213.05'
Set-TextValue ($ws.Range('E5')) '  +0.65%  '
Set-TextValue ($ws.Range('E6')) '  +0.17%  '
Set-TextValue ($ws.Range('D7')) '1.00'
Set-TextValue ($ws.Range('E7')) '  -0.79%  '
Set-TextValue ($ws.Range('D8')) '23.45'
Set-TextValue ($ws.Range('E8')) '  +6.45%  '
Set-TextValue ($ws.Range('E9')) '  +0.76%  '
Set-TextValue ($ws.Range('E10')) '  -0.17%  '
Set-TextValue ($ws.Range('E11')) '  +1.93%  '
Set-TextValue ($ws.Range('D12')) '1.802.24'
Set-TextValue ($ws.Range('E12')) '  +0.50%  '
Set-TextValue ($ws.Range('D13')) '1.579.80'
Set-TextValue ($ws.Range('E13')) '  +1.38%  '
Set-TextValue ($ws.Range('E14')) '  -0.79%  '
Set-TextValue ($ws.Range('E15')) '  +1.04%  '
Set-TextValue ($ws.Range('D16')) '28.130.81'
Set-TextValue ($ws.Range('E16')) '  +3.41%  '
Set-TextValue ($ws.Range('D17')) '63.61'
Set-TextValue ($ws.Range('E17')) '  +2.04%  '
Set-TextValue ($ws.Range('D18')) '229.73'
Set-TextValue ($ws.Range('E18')) '  +6.18%  '
Set-TextValue ($ws.Range('E19')) '  +0.43%  '
Set-TextValue ($ws.Range('E20')) '  +0.41%  '
Set-TextValue ($ws.Range('D21')) '1.00'
Set-TextValue ($ws.Range('E21')) '  -0.91%  '
Set-TextValue ($ws.Range('E22')) '  -0.65%  '
Set-TextValue ($ws.Range('D23')) '9.33'
Set-TextValue ($ws.Range('E23')) '  +0.98%  '
Set-TextValue ($ws.Range('E24')) '  -0.05%  '
Set-TextValue ($ws.Range('D25')) '152.37'
Set-TextValue ($ws.Range('E25')) '  -1.09%  '
Set-TextValue ($ws.Range('E26')) '  +0.91%  '
Set-TextValue ($ws.Range('E27')) '  -1.55%  '
Set-TextValue ($ws.Range('E28')) '  +0.12%  '
Set-TextValue ($ws.Range('D30')) '1.14'
Set-TextValue ($ws.Range('E30')) '  +0.13%  '
Set-TextValue ($ws.Range('E31')) '  +0.17%  '
Set-TextValue ($ws.Range('E32')) '  -0.82%  '
Set-TextValue ($ws.Range('D33')) '3.14'
Set-TextValue ($ws.Range('E33')) '  -1.12%  '
Set-TextValue ($ws.Range('D34')) '1.417.57'
Set-TextValue ($ws.Range('E34')) '  -2.48%  '
Set-TextValue ($ws.Range('E35')) '  -1.43%  '
Set-TextValue ($ws.Range('D36')) '1.05'
Set-TextValue ($ws.Range('E36')) '  -5.03%  '
Set-TextValue ($ws.Range('D37')) '2.33'
Set-TextValue ($ws.Range('E37')) '  -1.00%  '
Set-TextValue ($ws.Range('E38')) '  -0.20%  '
Set-TextValue ($ws.Range('D39')) '0.541'
Set-TextValue ($ws.Range('E39')) '  +1.01%  '
Set-TextValue ($ws.Range('E40')) '  +5.80%  '
Set-TextValue ($ws.Range('D41')) '0.808'
Set-TextValue ($ws.Range('E41')) '  -0.18%  '
Set-TextValue ($ws.Range('D42')) '1.00'
Set-TextValue ($ws.Range('E42')) '  -0.94%  '
Set-TextValue ($ws.Range('D43')) '5.63'
Set-TextValue ($ws.Range('E44')) '  -2.60%  '
Set-TextValue ($ws.Range('D45')) '1.81'
Set-TextValue ($ws.Range('E45')) '  +4.48%  '
Set-TextValue ($ws.Range('D46')) '63.73'
Set-TextValue ($ws.Range('E46')) '  -1.57%  '
Set-TextValue ($ws.Range('D47')) '1.714.90'
Set-TextValue ($ws.Range('E47')) '  +0.45%  '
Set-TextValue ($ws.Range('D48')) '87.11'
Set-TextValue ($ws.Range('E48')) '  +1.42%  '
Set-TextValue ($ws.Range('D49')) '0.0₆0105'
Set-TextValue ($ws.Range('E49')) '  +0.92%  '
Set-TextValue ($ws.Range('D50')) '0.0525'
Set-TextValue ($ws.Range('E50')) '  +0.96%  '
Set-TextValue ($ws.Range('E51')) '  -1.39%  '
